$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before H: old H (json-per-row formula) becomes I,
# old I (CONCAT aggregate formula) becomes J. Excel auto-shifts all the
# formula references (H:H -> I:I, H3 -> I3, etc.) when inserting this way.
$ws.Columns("H:H").Insert()

# --- Header row (row 1) ---
# (Order of assignment controls shared-string insertion order so the
# saved sharedStrings.xml table matches the original author's edit.)
$ws.Range("H1").Value = "order"
$ws.Range("G1").Value = "x"
$ws.Range("F1").Value = "y"

# --- "display" column E: "ture" -> "false" ---
# E2 is a literal text value; force text (not boolean) via quote-prefix.
$ws.Range("E2").Value = "'false"
# E3 references E2; E4:E11 reference the previous row (shared formula).
$ws.Range("E3").Formula = "=E2"
$ws.Range("E4:E11").Formula = "=E3"
# Re-stamp the quote-prefix text style (lost when the formulas were typed
# in) onto E3:E11, matching E2's style.
$ws.Range("E2").Copy()
$ws.Range("E3:E11").PasteSpecial(-4122)

# --- New "order" column H: sequential index 0..9 for rows 2..11 ---
for ($r = 2; $r -le 11; $r++) {
    $ws.Cells.Item($r, 8).Value = $r - 2
}

# --- Update the per-row JSON-builder formulas (col I) to include col H ---
$ws.Range("I2").Formula = '=CHAR(34)&A2&CHAR(34)&":{"&CHAR(34)&$B$1&CHAR(34)&":"&CHAR(34)&B2&CHAR(34)&","&CHAR(34)&$C$1&CHAR(34)&":"&CHAR(34)&C2&CHAR(34)&","&CHAR(34)&$D$1&CHAR(34)&":"&CHAR(34)&D2&CHAR(34)&","&CHAR(34)&$E$1&CHAR(34)&":"&CHAR(34)&E2&CHAR(34)&","&CHAR(34)&$F$1&CHAR(34)&":"&CHAR(34)&F2&CHAR(34)&","&CHAR(34)&$G$1&CHAR(34)&":"&CHAR(34)&G2&CHAR(34)&","&CHAR(34)&$H$1&CHAR(34)&":"&CHAR(34)&H2&CHAR(34)&"}"&IF(ISBLANK(A3),"",",")'

$ws.Range("I3").Formula = '=CHAR(34)&A3&CHAR(34)&":{"&CHAR(34)&$B$1&CHAR(34)&":"&CHAR(34)&B3&CHAR(34)&","&CHAR(34)&$C$1&CHAR(34)&":"&CHAR(34)&C3&CHAR(34)&","&CHAR(34)&$D$1&CHAR(34)&":"&CHAR(34)&D3&CHAR(34)&","&CHAR(34)&$E$1&CHAR(34)&":"&CHAR(34)&E3&CHAR(34)&","&CHAR(34)&$F$1&CHAR(34)&":"&CHAR(34)&F3&CHAR(34)&","&CHAR(34)&$G$1&CHAR(34)&":"&CHAR(34)&G3&CHAR(34)&","&CHAR(34)&$H$1&CHAR(34)&":"&CHAR(34)&H3&CHAR(34)&"}"&IF(ISBLANK(A4),"",",")'
$ws.Range("I3:I11").Formula = $ws.Range("I3").Formula

# --- Selection moves to J1 (mirrors the saved workbook view state) ---
$ws.Range("J1").Select()
